$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number-format on the Price/Volume columns before writing so
# Excel does not auto-coerce numeric-looking text (e.g. "1.001") into numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '30.646.75'
$ws.Range("E2").Value = '  +1.51%  '

$ws.Range("D3").Value = '1.861.98'
$ws.Range("E3").Value = '  +0.37%  '

$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = '235.57'
$ws.Range("E5").Value = '  +1.49%  '

$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.14%  '

$ws.Range("D7").Value = '0.4710'
$ws.Range("E7").Value = '  -0.33%  '

$ws.Range("D8").Value = '0.2761'
$ws.Range("E8").Value = '  +0.67%  '

$ws.Range("D9").Value = '0.06367'
$ws.Range("E9").Value = '  -0.66%  '

$ws.Range("D10").Value = '17.61'
$ws.Range("E10").Value = '  +8.52%  '

$ws.Range("D11").Value = '1.869.41'
$ws.Range("E11").Value = '  +1.02%  '

$ws.Range("D12").Value = '0.07461'
$ws.Range("E12").Value = '  +0.27%  '

$ws.Range("D13").Value = '5.126'
$ws.Range("E13").Value = '  +2.75%  '

$ws.Range("D14").Value = '85.02'
$ws.Range("E14").Value = '  -0.16%  '

$ws.Range("D15").Value = '0.6331'
$ws.Range("E15").Value = '  +0.45%  '

$ws.Range("D16").Value = '30.615.53'
$ws.Range("E16").Value = '  +1.53%  '

$ws.Range("D17").Value = '243.09'
$ws.Range("E17").Value = '  +5.16%  '

$ws.Range("D18").Value = '1.002'
$ws.Range("E18").Value = '  +0.11%  '

$ws.Range("D19").Value = '12.88'
$ws.Range("E19").Value = '  +1.40%  '

$ws.Range("D20").Value = '0.000007372'
$ws.Range("E20").Value = '  +0.64%  '

$ws.Range("E21").Value = '  -0.02%  '

$ws.Range("D22").Value = '5.018'
$ws.Range("E22").Value = '  -0.25%  '

$ws.Range("D23").Value = '6.089'
$ws.Range("E23").Value = '  +1.50%  '

$ws.Range("D24").Value = '9.381'
$ws.Range("E24").Value = '  +1.52%  '

$ws.Range("D25").Value = '164.76'
$ws.Range("E25").Value = '  -0.51%  '

$ws.Range("D26").Value = '18.18'
$ws.Range("E26").Value = '  +2.10%  '

$ws.Range("D27").Value = '1.893'
$ws.Range("E27").Value = '  +0.47%  '

$ws.Range("D28").Value = '0.1016'
$ws.Range("E28").Value = '  -0.45%  '

$ws.Range("D29").Value = '1.385'
$ws.Range("E29").Value = '  +0.22%  '

$ws.Range("D30").Value = '4.068'
$ws.Range("E30").Value = '  -1.34%  '

$ws.Range("D31").Value = '3.867'
$ws.Range("E31").Value = '  -1.24%  '

$ws.Range("D32").Value = '0.04954'
$ws.Range("E32").Value = '  +1.08%  '

$ws.Range("D33").Value = '1.151'
$ws.Range("E33").Value = '  +0.51%  '

$ws.Range("D34").Value = '0.7063'
$ws.Range("E34").Value = '  -1.82%  '

$ws.Range("D35").Value = '2.716'
$ws.Range("E35").Value = '  +0.98%  '

$ws.Range("D36").Value = '0.01910'
$ws.Range("E36").Value = '  +1.06%  '

$ws.Range("D37").Value = '2.683'
$ws.Range("E37").Value = '  +1.77%  '

$ws.Range("D38").Value = '0.8815'
$ws.Range("E38").Value = '  -2.38%  '

$ws.Range("E39").Value = '  +2.32%  '

$ws.Range("D40").Value = '105.28'
$ws.Range("E40").Value = '  -0.34%  '

$ws.Range("D41").Value = '1.002'
$ws.Range("E41").Value = '  +0.31%  '

$ws.Range("D42").Value = '5.573'
$ws.Range("E42").Value = '  +0.79%  '

$ws.Range("D43").Value = '0.4094'
$ws.Range("E43").Value = '  -0.04%  '

$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = '7.265'
$ws.Range("E44").Value = '  +2.64%  '

$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '64.40'
$ws.Range("E45").Value = '  +5.42%  '

$ws.Range("D46").Value = '0.1218'
$ws.Range("E46").Value = '  +1.86%  '

$ws.Range("D47").Value = '33.66'
$ws.Range("E47").Value = '  +1.42%  '

$ws.Range("D48").Value = '8.627'
$ws.Range("E48").Value = '  -0.99%  '

$ws.Range("D49").Value = '0.05558'
$ws.Range("E49").Value = '  -0.28%  '

$ws.Range("D50").Value = '1.374'
$ws.Range("E50").Value = '  -1.92%  '

$ws.Range("D51").Value = '0.3704'
$ws.Range("E51").Value = '  +0.44%  '

# Restore the default (Normal) cell style so no stray number-format
# artifacts are left behind from the text-coercion step above.
$dataRange.Style = "Normal"
